# Apply the cryptos-list data refresh described in the commit diff
# ("Updated cryptos list on Sun Sep 22 20:27:46 UTC 2024 with GitHub Actions").
#
# The Price column (D) stores values as text, not numbers -- e.g. "63.152.07" or
# "0.0000145" are not valid Excel numbers, and even genuinely numeric-looking prices
# such as "592.34" are kept as text in this sheet. Assigning a numeric-looking string
# straight to Range.Value would make Excel auto-convert it to a real number and change
# the cell from a string cell into a numeric cell. To avoid that, any refreshed Price
# value that would parse as a plain number is written with a leading apostrophe, which
# is the standard Excel convention for forcing text entry; Excel strips the apostrophe
# and stores the digits as text, exactly like the rest of the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.235.01'
$ws.Range("E2").Value = '  +0.06%  '
$ws.Range("D3").Value = '2.578.47'
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''592.34'
$ws.Range("E5").Value = '  +1.21%  '
$ws.Range("D6").Value = '''144.89'
$ws.Range("E6").Value = '  -1.78%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '''0.591'
$ws.Range("E8").Value = '  -1.80%  '
$ws.Range("E9").Value = '  -2.05%  '
$ws.Range("D10").Value = '''5.60'
$ws.Range("E10").Value = '  -0.69%  '
$ws.Range("E11").Value = '  -0.12%  '
$ws.Range("D12").Value = '''0.352'
$ws.Range("E12").Value = '  -1.46%  '
$ws.Range("D13").Value = '''27.19'
$ws.Range("E13").Value = '  -1.10%  '
$ws.Range("D14").Value = '3.042.34'
$ws.Range("E14").Value = '  +0.42%  '
$ws.Range("D15").Value = '63.140.24'
$ws.Range("E15").Value = '  +0.01%  '
$ws.Range("D16").Value = '''0.0000146'
$ws.Range("E16").Value = '  -1.06%  '
$ws.Range("D17").Value = '2.577.10'
$ws.Range("E17").Value = '  +0.36%  '
$ws.Range("D18").Value = '''11.11'
$ws.Range("E18").Value = '  -2.32%  '
$ws.Range("D19").Value = '''342.46'
$ws.Range("E19").Value = '  -0.41%  '
$ws.Range("D20").Value = '''4.35'
$ws.Range("E23").Value = '  +3.61%  '
$ws.Range("D24").Value = '''67.88'
$ws.Range("E24").Value = '  +1.49%  '
$ws.Range("D25").Value = '''1.60'
$ws.Range("E25").Value = '  +8.44%  '
$ws.Range("E26").Value = '  -0.39%  '
$ws.Range("D27").Value = '''0.166'
$ws.Range("E27").Value = '  -3.04%  '
$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").Value = '''0.999'
$ws.Range("E28").Value = '  -0.14%  '
$ws.Range("B29").Value = 'Aptos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D29").Value = '''7.94'
$ws.Range("E29").Value = '  -2.50%  '
$ws.Range("E30").Value = '  -2.82%  '
$ws.Range("E31").Value = '  -2.25%  '
$ws.Range("D32").Value = '''468.85'
$ws.Range("E32").Value = '  +0.61%  '
$ws.Range("D33").Value = '0.0₃0802'
$ws.Range("E33").Value = '  -3.16%  '
$ws.Range("E34").Value = '  +3.24%  '
$ws.Range("D35").Value = '''176.58'
$ws.Range("E35").Value = '  +0.25%  '
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("D37").Value = '''0.396'
$ws.Range("E37").Value = '  -3.08%  '
$ws.Range("D38").Value = '''18.89'
$ws.Range("E38").Value = '  -1.84%  '
$ws.Range("E39").Value = '  +0.44%  '
$ws.Range("D40").Value = '''0.999'
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("E41").Value = '  -3.29%  '
$ws.Range("D42").Value = '''158.55'
$ws.Range("D43").Value = '''40.03'
$ws.Range("E43").Value = '  +1.06%  '
$ws.Range("D44").Value = '''3.71'
$ws.Range("E44").Value = '  -3.00%  '
$ws.Range("D45").Value = '''21.47'
$ws.Range("E45").Value = '  +2.27%  '
$ws.Range("D46").Value = '''0.635'
$ws.Range("E46").Value = '  +3.63%  '
$ws.Range("D47").Value = '''0.0539'
$ws.Range("E47").Value = '  -1.67%  '
$ws.Range("E48").Value = '  -1.64%  '
$ws.Range("D49").Value = '''0.0238'
$ws.Range("E49").Value = '  -0.71%  '
$ws.Range("D50").Value = '''18.19'
$ws.Range("E50").Value = '  -1.93%  '
